$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 391.3
$ws.Range("I5").Value = 398.66666
$ws.Range("J5").Value = 325
$ws.Range("K5").Value = 398.66666
$ws.Range("L5").Value = 325
$ws.Range("M5").Value = -283.66666
$ws.Range("N5").Value = -555
$ws.Range("H33").Value = 561.36365
$ws.Range("I33").Value = 444.6
$ws.Range("J33").Value = 658.6667
$ws.Range("K33").Value = 444.6
$ws.Range("L33").Value = 658.6667
$ws.Range("M33").Value = -215.6
$ws.Range("N33").Value = -1116.6667
$ws.Range("H48").Value = 999.2857
$ws.Range("I48").Value = 999.2857
$ws.Range("K48").Value = 2997.8571
$ws.Range("M48").Value = -2705.8571
$ws.Range("H56").Value = 999.2857
$ws.Range("I56").Value = 999.2857
$ws.Range("K56").Value = 2997.8571
$ws.Range("M56").Value = -2463.8571
$ws.Range("H74").Value = 9823
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = $null
$ws.Range("H76").Value = 7410.7144
$ws.Range("I76").Value = 5937.5
$ws.Range("K76").Value = 5937.5
$ws.Range("M76").Value = -5622.5
$ws.Range("H77").Value = 9823
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = $null
$ws.Range("H79").Value = 7410.7144
$ws.Range("I79").Value = 5937.5
$ws.Range("K79").Value = 5937.5
$ws.Range("M79").Value = -4845.5
$ws.Range("H82").Value = 2905.125
$ws.Range("I82").Value = 2905.125
$ws.Range("K82").Value = 8715.375
$ws.Range("M82").Value = -8309.375
$ws.Range("H85").Value = 2905.125
$ws.Range("I85").Value = 2905.125
$ws.Range("K85").Value = 8715.375
$ws.Range("M85").Value = -7311.375
$ws.Range("H103").Value = 825.5
$ws.Range("I103").Value = 538.25
$ws.Range("K103").Value = 1614.75
$ws.Range("M103").Value = -1028.75
$ws.Range("H112").Value = 1731.0769
$ws.Range("J112").Value = 1888.091
$ws.Range("L112").Value = 5664.272999999999
$ws.Range("N112").Value = -7880.272999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3792069.5
$ws.Range("I32").Value = 4117915.5
$ws.Range("K32").Value = 4117915.5
$ws.Range("M32").Value = -4117628.5
$ws.Range("H45").Value = 5287.25
$ws.Range("J45").Value = 5309.1665
$ws.Range("L45").Value = 5309.1665
$ws.Range("N45").Value = -6063.1665
$ws.Range("H61").Value = 1974060.2
$ws.Range("I61").Value = 2096376.5
$ws.Range("J61").Value = 17000
$ws.Range("K61").Value = 2096376.5
$ws.Range("L61").Value = 17000
$ws.Range("M61").Value = -2096164.5
$ws.Range("N61").Value = -17424
$ws.Range("H63").Value = 5593.579
$ws.Range("I63").Value = 2282.1667
$ws.Range("J63").Value = 7121.923
$ws.Range("K63").Value = 2282.1667
$ws.Range("L63").Value = 7121.923
$ws.Range("M63").Value = -1596.1667
$ws.Range("N63").Value = -8493.922999999999
$ws.Range("H66").Value = 5593.579
$ws.Range("I66").Value = 2282.1667
$ws.Range("J66").Value = 7121.923
$ws.Range("K66").Value = 11410.8335
$ws.Range("L66").Value = 35609.615
$ws.Range("M66").Value = -7978.833500000001
$ws.Range("N66").Value = -42473.615
$ws.Range("H74").Value = 1955718.8
$ws.Range("I74").Value = 2606285
$ws.Range("K74").Value = 2606285
$ws.Range("M74").Value = -2605411
$ws.Range("H77").Value = 1955718.8
$ws.Range("I77").Value = 2606285
$ws.Range("K77").Value = 13031425
$ws.Range("M77").Value = -13027057
$ws.Range("H97").Value = 1263.6923
$ws.Range("I97").Value = 1302.4166
$ws.Range("K97").Value = 1302.4166
$ws.Range("M97").Value = -806.4166
$ws.Range("H122").Value = 3199.5
$ws.Range("I122").Value = 1910.2222
$ws.Range("J122").Value = 4857.143
$ws.Range("K122").Value = 5730.6666
$ws.Range("L122").Value = 14571.429
$ws.Range("M122").Value = -3280.6666
$ws.Range("N122").Value = -19471.429
$ws.Range("H132").Value = 687909.6
$ws.Range("I132").Value = 778070.9399999999
$ws.Range("J132").Value = 11699.75
$ws.Range("K132").Value = 2334212.82
$ws.Range("L132").Value = 35099.25
$ws.Range("M132").Value = -2331682.82
$ws.Range("N132").Value = -40159.25
$ws.Range("H136").Value = 1974060.2
$ws.Range("I136").Value = 2096376.5
$ws.Range("J136").Value = 17000
$ws.Range("K136").Value = 6289129.5
$ws.Range("L136").Value = 51000
$ws.Range("M136").Value = -6286579.5
$ws.Range("N136").Value = -56100
$ws.Range("H139").Value = 175000
$ws.Range("J139").Value = 175000
$ws.Range("L139").Value = 175000
$ws.Range("N139").Value = -185280

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1881.7391
$ws.Range("I20").Value = 1854.9412
$ws.Range("K20").Value = 1854.9412
$ws.Range("M20").Value = -1607.9412
$ws.Range("H58").Value = 90000
$ws.Range("J58").Value = 90000
$ws.Range("L58").Value = 90000
$ws.Range("N58").Value = -90588
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = $null
$ws.Range("M60").Value = $null
$ws.Range("N60").Value = 0
$ws.Range("H94").Value = 2250.2856
$ws.Range("I94").Value = 950.8
$ws.Range("K94").Value = 950.8
$ws.Range("M94").Value = -499.8
$ws.Range("H105").Value = 2175.4285
$ws.Range("I105").Value = 2068.3333
$ws.Range("K105").Value = 2068.3333
$ws.Range("M105").Value = -321.3332999999998
$ws.Range("H123").Value = 40000
$ws.Range("I123").Value = 40000
$ws.Range("K123").Value = 40000
$ws.Range("M123").Value = -35100
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = $null
$ws.Range("N132").Value = 0
$ws.Range("H134").Value = 755880.9
$ws.Range("I134").Value = 921013.9
$ws.Range("J134").Value = 12782.333
$ws.Range("K134").Value = 2763041.7
$ws.Range("L134").Value = 38346.999
$ws.Range("M134").Value = -2760506.7
$ws.Range("N134").Value = -43416.999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8852.036
$ws.Range("I31").Value = 4042.7778
$ws.Range("J31").Value = 11130.105
$ws.Range("K31").Value = 4042.7778
$ws.Range("L31").Value = 11130.105
$ws.Range("M31").Value = -3747.7778
$ws.Range("N31").Value = -11720.105
$ws.Range("H34").Value = 8852.036
$ws.Range("I34").Value = 4042.7778
$ws.Range("J34").Value = 11130.105
$ws.Range("K34").Value = 4042.7778
$ws.Range("L34").Value = 11130.105
$ws.Range("M34").Value = -3840.7778
$ws.Range("N34").Value = -11534.105
$ws.Range("H58").Value = 1238583.5
$ws.Range("I58").Value = 1765406.4
$ws.Range("J58").Value = 9330
$ws.Range("K58").Value = 1765406.4
$ws.Range("L58").Value = 9330
$ws.Range("M58").Value = -1765203.4
$ws.Range("N58").Value = -9736
$ws.Range("H99").Value = 5412.5
$ws.Range("I99").Value = 4143.75
$ws.Range("K99").Value = 4143.75
$ws.Range("M99").Value = -2645.75
$ws.Range("H122").Value = 6854.2856
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 7330
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 21990
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -26890
$ws.Range("H123").Value = 299999
$ws.Range("J123").Value = 299999
$ws.Range("L123").Value = 299999
$ws.Range("N123").Value = -309799
$ws.Range("H126").Value = 5412.5
$ws.Range("I126").Value = 4143.75
$ws.Range("K126").Value = 12431.25
$ws.Range("M126").Value = -9961.25
$ws.Range("H132").Value = 71507.8
$ws.Range("I132").Value = 71507.8
$ws.Range("K132").Value = 214523.4
$ws.Range("M132").Value = -211993.4
$ws.Range("H134").Value = 1927.8889
$ws.Range("I134").Value = 1993.875
$ws.Range("J134").Value = 1400
$ws.Range("K134").Value = 5981.625
$ws.Range("L134").Value = 4200
$ws.Range("M134").Value = -3446.625
$ws.Range("N134").Value = -9270
$ws.Range("H136").Value = 1238583.5
$ws.Range("I136").Value = 1765406.4
$ws.Range("J136").Value = 9330
$ws.Range("K136").Value = 5296219.199999999
$ws.Range("L136").Value = 27990
$ws.Range("M136").Value = -5293669.199999999
$ws.Range("N136").Value = -33090

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 162.66667
$ws.Range("I41").Value = 100
$ws.Range("K41").Value = 300
$ws.Range("M41").Value = 38
$ws.Range("H140").Value = 3129.7368
$ws.Range("I140").Value = 2748.0557
$ws.Range("K140").Value = 8244.167099999999
$ws.Range("M140").Value = -3064.167099999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 5190000
$ws.Range("I14").Value = 400000
$ws.Range("J14").Value = 9980000
$ws.Range("K14").Value = 400000
$ws.Range("L14").Value = 9980000
$ws.Range("M14").Value = -399832
$ws.Range("N14").Value = -9980336
$ws.Range("H102").Value = 2234.3333
$ws.Range("I102").Value = 2234.3333
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2234.3333
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = $null
$ws.Range("N102").Value = -612.3332999999998
$ws.Range("H122").Value = 89900
$ws.Range("I122").Value = 29800
$ws.Range("K122").Value = 89400
$ws.Range("M122").Value = -86950
$ws.Range("H126").Value = 837668.8
$ws.Range("I126").Value = 1668613.9
$ws.Range("J126").Value = 6723.7
$ws.Range("K126").Value = 5005841.699999999
$ws.Range("L126").Value = 20171.1
$ws.Range("M126").Value = -5003371.699999999
$ws.Range("N126").Value = -25111.1
$ws.Range("H132").Value = 2093.125
$ws.Range("I132").Value = 1229.909
$ws.Range("K132").Value = 3689.727
$ws.Range("M132").Value = -1159.727
$ws.Range("H136").Value = 21001.584
$ws.Range("I136").Value = 12325
$ws.Range("J136").Value = 21790.363
$ws.Range("K136").Value = 36975
$ws.Range("L136").Value = 65371.08900000001
$ws.Range("M136").Value = -34425
$ws.Range("N136").Value = -70471.08900000001
$ws.Range("H137").Value = 72694.625
$ws.Range("J137").Value = 110780
$ws.Range("L137").Value = 110780
$ws.Range("N137").Value = -120980

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I12").Value = 4500
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 4500
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = $null
$ws.Range("N12").Value = -4330
$ws.Range("H82").Value = 2487.6667
$ws.Range("I82").Value = 2240.625
$ws.Range("J82").Value = 2770
$ws.Range("K82").Value = 2240.625
$ws.Range("L82").Value = 2770
$ws.Range("M82").Value = -1879.625
$ws.Range("N82").Value = -3492
$ws.Range("H85").Value = 2487.6667
$ws.Range("I85").Value = 2240.625
$ws.Range("J85").Value = 2770
$ws.Range("K85").Value = 2240.625
$ws.Range("L85").Value = 2770
$ws.Range("M85").Value = -992.625
$ws.Range("N85").Value = -5266
$ws.Range("H122").Value = 4631.45
$ws.Range("I122").Value = 4424.3076
$ws.Range("J122").Value = 5016.143
$ws.Range("K122").Value = 13272.9228
$ws.Range("L122").Value = 15048.429
$ws.Range("M122").Value = -10822.9228
$ws.Range("N122").Value = -19948.429
$ws.Range("H132").Value = 3352.3171
$ws.Range("I132").Value = 2998.4285
$ws.Range("K132").Value = 8995.2855
$ws.Range("M132").Value = -6465.2855
$ws.Range("H136").Value = 5418
$ws.Range("I136").Value = 4959.8
$ws.Range("K136").Value = 14879.4
$ws.Range("M136").Value = -12329.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 20013
$ws.Range("I32").Value = 20013
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 20013
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = $null
$ws.Range("N32").Value = -19696
$ws.Range("H122").Value = 3981.0908
$ws.Range("I122").Value = 3482
$ws.Range("J122").Value = 4580
$ws.Range("K122").Value = 10446
$ws.Range("L122").Value = 13740
$ws.Range("M122").Value = -7996
$ws.Range("N122").Value = -18640
$ws.Range("H126").Value = 5225.5454
$ws.Range("I126").Value = 4435.875
$ws.Range("J126").Value = 7331.3335
$ws.Range("K126").Value = 13307.625
$ws.Range("L126").Value = 21994.0005
$ws.Range("M126").Value = -10837.625
$ws.Range("N126").Value = -26934.0005
$ws.Range("H132").Value = 9266781
$ws.Range("I132").Value = 16674516
$ws.Range("J132").Value = 7112.375
$ws.Range("K132").Value = 50023548
$ws.Range("L132").Value = 21337.125
$ws.Range("M132").Value = -50021018
$ws.Range("N132").Value = -26397.125

Write-Host "Applied all changes"